$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.696.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.98%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.160.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.80%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +16.32%  "

# Row 9
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("E10").Value = "  +5.40%  "

# Row 11
$ws.Range("E11").Value = "  +3.90%  "

# Row 12
$ws.Range("E12").Value = "  +3.73%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.703.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.09%  "

# Row 15
$ws.Range("E15").Value = "  +4.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.740.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "

# Row 17
$ws.Range("E17").Value = "  +3.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.135.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.90%  "

# Row 20
$ws.Range("E20").Value = "  -0.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.90%  "

# Row 22
$ws.Range("E22").Value = "  +2.15%  "

# Row 23
$ws.Range("E23").Value = "  +0.14%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.533"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.12%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.81%  "

# Row 26
$ws.Range("E26").Value = "  +0.56%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.94%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0868"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.16%  "

# Row 31
$ws.Range("E31").Value = "  +0.06%  "

# Row 32
$ws.Range("E32").Value = "  -0.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.59%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.18%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "

# Row 37
$ws.Range("E37").Value = "  +4.68%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.67%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.15%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0689"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.37%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.647.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.37%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.722"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.58%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0289"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.92%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.200.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "

# Row 48
$ws.Range("E48").Value = "  +14.29%  "

# Row 49
$ws.Range("E49").Value = "  +1.86%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.978"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.76%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "
